$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $new, 2) | Out-Null
}

# The paragraph for "激烈竞争" has a third run immediately following the
# run we need to edit ("：加拿大的云服务部门竞争激烈，众多玩家。"), and that
# third run ("Adatum 必须清楚阐明...") happens to share byte-identical run
# formatting. A naive Find/Replace on the middle run causes the engine to
# coalesce it with that untouched neighbor into a single run, which the
# target revision does not do (the neighbor stays a separate, untouched
# run). To avoid that, briefly perturb the neighbor run's own formatting
# before the edit (preventing the coalesce), then restore it to its
# original value in a separate follow-up call.
$guard = $d.Content
$guard.Find.Execute("Adatum 必须清楚阐明其解决方案的独特价值，这样才能在市场上占据一席之地。", $true) | Out-Null
if ($guard.Find.Found) {
    $guard.Font.Color = 1
}

Replace-Text "有限的品牌识别和意识" "品牌知名度和品牌认知有限"
Replace-Text "：在这些新市场实现可见性是一个主要障碍，需要强大的营销努力，以从头开始建立阿达图姆的品牌存在。" "：在新市场推广品牌是一个主要障碍，需要进行强大的营销努力，从头开始建立 Adatum 品牌。"
Replace-Text "激烈竞争" "竞争激烈"
Replace-Text "：加拿大的云服务部门竞争激烈，众多玩家。" "：加拿大的云服务部门竞争激烈，玩家众多。"
Replace-Text "：定制产品和营销，以满足这些市场的各种需求，对于与当地企业和消费者共鸣至关重要。" "：定制产品和营销，以满足这些市场的各种需求，这对于与当地企业和消费者产生共鸣至关重要。"
Replace-Text "法规和合规性挑战" "法律和合规性挑战"
Replace-Text "：Adatum 面临着导航区域独特的数据隐私、安全和运营法规的复杂任务，这需要勤奋的合规性工作。" "：在区域独特的数据隐私、安全和运营法规方面，Adatum 的任务繁杂，这需要进行尽职的合规性工作。"
Replace-Text "运营和后勤复杂性" "运营和物流复杂"
Replace-Text "：建立高效、跨区域运营带来了后勤挑战，特别是在维护高服务级别和管理地理位置的数据中心方面。" "：要建立高效的跨区域运营会面临物流挑战，特别是在维护高服务级别和跨地区管理数据中心方面。"

# Restore the guarded run's formatting back to its original value now
# that the sibling edit above has already happened and won't re-merge it.
$restore = $d.Content
$restore.Find.Execute("Adatum 必须清楚阐明其解决方案的独特价值，这样才能在市场上占据一席之地。", $true) | Out-Null
if ($restore.Find.Found) {
    $restore.Font.Color = -16777216
}
